# Fix the algorithm/conditions on filtering the status of candidates
# Re-sorts the Active Candidates list by Job ID and adds 3 new candidate rows
# (Safe Security x2, Axion Ray - Matt Henry), extending the range to A1:F27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 594
$ws.Cells.Item(2, 2).Value = "Antithesis"
$ws.Cells.Item(2, 3).Value = "Enterprise Account Executive - East"
$ws.Cells.Item(2, 4).Value = "Patrick Racy"
$ws.Cells.Item(2, 5).Value = "2nd Interview"
$ws.Cells.Item(2, 6).Value = 45984

$ws.Cells.Item(3, 1).Value = 594
$ws.Cells.Item(3, 2).Value = "Antithesis"
$ws.Cells.Item(3, 3).Value = "Enterprise Account Executive - East"
$ws.Cells.Item(3, 4).Value = "Ryan Lallier"
$ws.Cells.Item(3, 5).Value = "1st Interview"
$ws.Cells.Item(3, 6).Value = 45987

$ws.Cells.Item(4, 1).Value = 665
$ws.Cells.Item(4, 2).Value = "Safe Security"
$ws.Cells.Item(4, 3).Value = "Enterprise AE UK"
$ws.Cells.Item(4, 4).Value = "Iain Airey"
$ws.Cells.Item(4, 5).Value = "3rd Interview"
$ws.Cells.Item(4, 6).Value = 45975

$ws.Cells.Item(5, 1).Value = 665
$ws.Cells.Item(5, 2).Value = "Safe Security"
$ws.Cells.Item(5, 3).Value = "Enterprise AE UK"
$ws.Cells.Item(5, 4).Value = "Connor Newby"
$ws.Cells.Item(5, 5).Value = "4th Interview"
$ws.Cells.Item(5, 6).Value = 45968

$ws.Cells.Item(6, 1).Value = 667
$ws.Cells.Item(6, 2).Value = "Antithesis"
$ws.Cells.Item(6, 3).Value = "BDR"
$ws.Cells.Item(6, 4).Value = "Daniel Jones"
$ws.Cells.Item(6, 5).Value = "CV Sent"
$ws.Cells.Item(6, 6).Value = 45987

$ws.Cells.Item(7, 1).Value = 667
$ws.Cells.Item(7, 2).Value = "Antithesis"
$ws.Cells.Item(7, 3).Value = "BDR"
$ws.Cells.Item(7, 4).Value = "J. Donahoe"
$ws.Cells.Item(7, 5).Value = "3rd Interview"
$ws.Cells.Item(7, 6).Value = 45989

$ws.Cells.Item(8, 1).Value = 667
$ws.Cells.Item(8, 2).Value = "Antithesis"
$ws.Cells.Item(8, 3).Value = "BDR"
$ws.Cells.Item(8, 4).Value = "Omar Mohamed"
$ws.Cells.Item(8, 5).Value = "2nd Interview"
$ws.Cells.Item(8, 6).Value = 45989

$ws.Cells.Item(9, 1).Value = 680
$ws.Cells.Item(9, 2).Value = "Oscilar"
$ws.Cells.Item(9, 3).Value = "Sales Engineer"
$ws.Cells.Item(9, 4).Value = "Kyle Martinez"
$ws.Cells.Item(9, 5).Value = "2nd Interview"
$ws.Cells.Item(9, 6).Value = 45987

$ws.Cells.Item(10, 1).Value = 680
$ws.Cells.Item(10, 2).Value = "Oscilar"
$ws.Cells.Item(10, 3).Value = "Sales Engineer"
$ws.Cells.Item(10, 4).Value = "JOHN FROST"
$ws.Cells.Item(10, 5).Value = "2nd Interview"
$ws.Cells.Item(10, 6).Value = 45989

$ws.Cells.Item(11, 1).Value = 680
$ws.Cells.Item(11, 2).Value = "Oscilar"
$ws.Cells.Item(11, 3).Value = "Sales Engineer"
$ws.Cells.Item(11, 4).Value = "Andrew Barnett"
$ws.Cells.Item(11, 5).Value = "2nd Interview"
$ws.Cells.Item(11, 6).Value = 45984

$ws.Cells.Item(12, 1).Value = 680
$ws.Cells.Item(12, 2).Value = "Oscilar"
$ws.Cells.Item(12, 3).Value = "Sales Engineer"
$ws.Cells.Item(12, 4).Value = "Jonathan Valand"
$ws.Cells.Item(12, 5).Value = "2nd Interview"
$ws.Cells.Item(12, 6).Value = 45989

$ws.Cells.Item(13, 1).Value = 704
$ws.Cells.Item(13, 2).Value = "Dash0"
$ws.Cells.Item(13, 3).Value = "SDR Leader NYC"
$ws.Cells.Item(13, 4).Value = "Michael Vargas"
$ws.Cells.Item(13, 5).Value = "CV Sent"
$ws.Cells.Item(13, 6).Value = 45965

$ws.Cells.Item(14, 1).Value = 731
$ws.Cells.Item(14, 2).Value = "Oscilar"
$ws.Cells.Item(14, 3).Value = "Enterprise AE x5"
$ws.Cells.Item(14, 4).Value = "Jason Fulwood"
$ws.Cells.Item(14, 5).Value = "4th Interview"
$ws.Cells.Item(14, 6).Value = 45984

$ws.Cells.Item(15, 1).Value = 739
$ws.Cells.Item(15, 2).Value = "Axion Ray"
$ws.Cells.Item(15, 3).Value = "BDR Manager (North East)"
$ws.Cells.Item(15, 4).Value = "Tae Kim"
$ws.Cells.Item(15, 5).Value = "CV Sent"
$ws.Cells.Item(15, 6).Value = 45987

$ws.Cells.Item(16, 1).Value = 739
$ws.Cells.Item(16, 2).Value = "Axion Ray"
$ws.Cells.Item(16, 3).Value = "BDR Manager (North East)"
$ws.Cells.Item(16, 4).Value = "Matt Henry"
$ws.Cells.Item(16, 5).Value = "4th Interview"
$ws.Cells.Item(16, 6).Value = 45980

$ws.Cells.Item(17, 1).Value = 744
$ws.Cells.Item(17, 2).Value = "Synthflow AI"
$ws.Cells.Item(17, 3).Value = "Sales Engineers"
$ws.Cells.Item(17, 4).Value = "Nathan Watson"
$ws.Cells.Item(17, 5).Value = "CV Sent"
$ws.Cells.Item(17, 6).Value = 45980

$ws.Cells.Item(18, 1).Value = 744
$ws.Cells.Item(18, 2).Value = "Synthflow AI"
$ws.Cells.Item(18, 3).Value = "Sales Engineers"
$ws.Cells.Item(18, 4).Value = "Kenneth Winfield"
$ws.Cells.Item(18, 5).Value = "CV Sent"
$ws.Cells.Item(18, 6).Value = 45980

$ws.Cells.Item(19, 1).Value = 744
$ws.Cells.Item(19, 2).Value = "Synthflow AI"
$ws.Cells.Item(19, 3).Value = "Sales Engineers"
$ws.Cells.Item(19, 4).Value = "RICARDO DE BIASE"
$ws.Cells.Item(19, 5).Value = "2nd Interview"
$ws.Cells.Item(19, 6).Value = 45982

$ws.Cells.Item(20, 1).Value = 760
$ws.Cells.Item(20, 2).Value = "Impala"
$ws.Cells.Item(20, 3).Value = "Head of Sales (NA)"
$ws.Cells.Item(20, 4).Value = "Patrick Racy"
$ws.Cells.Item(20, 5).Value = "3rd Interview"
$ws.Cells.Item(20, 6).Value = 45989

$ws.Cells.Item(21, 1).Value = 807
$ws.Cells.Item(21, 2).Value = "Oscilar"
$ws.Cells.Item(21, 3).Value = "SE Leader"
$ws.Cells.Item(21, 4).Value = "Andrew Birnbaum"
$ws.Cells.Item(21, 5).Value = "1st Interview"
$ws.Cells.Item(21, 6).Value = 45982

$ws.Cells.Item(22, 1).Value = 807
$ws.Cells.Item(22, 2).Value = "Oscilar"
$ws.Cells.Item(22, 3).Value = "SE Leader"
$ws.Cells.Item(22, 4).Value = "Craig Fetterman"
$ws.Cells.Item(22, 5).Value = "1st Interview"
$ws.Cells.Item(22, 6).Value = 45982

$ws.Cells.Item(23, 1).Value = 807
$ws.Cells.Item(23, 2).Value = "Oscilar"
$ws.Cells.Item(23, 3).Value = "SE Leader"
$ws.Cells.Item(23, 4).Value = "David Hoenig"
$ws.Cells.Item(23, 5).Value = "1st Interview"
$ws.Cells.Item(23, 6).Value = 45987

$ws.Cells.Item(24, 1).Value = 807
$ws.Cells.Item(24, 2).Value = "Oscilar"
$ws.Cells.Item(24, 3).Value = "SE Leader"
$ws.Cells.Item(24, 4).Value = "Ray Mi"
$ws.Cells.Item(24, 5).Value = "3rd Interview"
$ws.Cells.Item(24, 6).Value = 45982

$ws.Cells.Item(25, 1).Value = 820
$ws.Cells.Item(25, 2).Value = "Silverfort"
$ws.Cells.Item(25, 3).Value = "Nordics RSM"
$ws.Cells.Item(25, 4).Value = "Jesper Damm-Skogh"
$ws.Cells.Item(25, 5).Value = "1st Interview"
$ws.Cells.Item(25, 6).Value = 45991
$ws.Cells.Item(25, 6).NumberFormat = $ws.Cells.Item(24, 6).NumberFormat

$ws.Cells.Item(26, 1).Value = 820
$ws.Cells.Item(26, 2).Value = "Silverfort"
$ws.Cells.Item(26, 3).Value = "Nordics RSM"
$ws.Cells.Item(26, 4).Value = "Marc Solis"
$ws.Cells.Item(26, 5).Value = "1st Interview"
$ws.Cells.Item(26, 6).Value = 45991
$ws.Cells.Item(26, 6).NumberFormat = $ws.Cells.Item(24, 6).NumberFormat

$ws.Cells.Item(27, 1).Value = 833
$ws.Cells.Item(27, 2).Value = "Blockaid"
$ws.Cells.Item(27, 3).Value = "SDR Manager"
$ws.Cells.Item(27, 4).Value = "Tae Kim"
$ws.Cells.Item(27, 5).Value = "CV Sent"
$ws.Cells.Item(27, 6).Value = 45987
$ws.Cells.Item(27, 6).NumberFormat = $ws.Cells.Item(24, 6).NumberFormat
